$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report issue number and week-of dates) ---
$ws.Range("A8").Value = "Volume 30   Number  24"
$ws.Range("C9").Value = "Report Covering the Week  6/12/2023  Through  6/18/2023"

# --- Type/style changes: text<->number conversions ---
# Use stable, untouched donor cells (rows 36+) to copy over exact cell style (format),
# then set the correct final value. For "0" text cells we copy directly from C14
# (which stays text "0" before and after this edit) so no value write is needed.

# F14: was number 2 -> becomes text "0" (style matches C14/D14)
$ws.Range("C14").Copy($ws.Range("F14"))

# C26: was number 1 -> becomes text "0"
$ws.Range("C14").Copy($ws.Range("C26"))

# F28: was number 2 -> becomes text "0"
$ws.Range("C14").Copy($ws.Range("F28"))

# F29: was number 2 -> becomes text "0"
$ws.Range("C14").Copy($ws.Range("F29"))

# D22: was text "0" -> becomes number 2
$ws.Range("C36").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 2

# E22: was text "***.*" -> becomes number -100
$ws.Range("K36").Copy($ws.Range("E22"))
$ws.Range("E22").Value = -100

# G22: was text "0" -> becomes number 2
$ws.Range("C36").Copy($ws.Range("G22"))
$ws.Range("G22").Value = 2

# H22: was text "***.*" -> becomes number -50
$ws.Range("K36").Copy($ws.Range("H22"))
$ws.Range("H22").Value = -50

# D28: was text "0" -> becomes number 5
$ws.Range("C36").Copy($ws.Range("D28"))
$ws.Range("D28").Value = 5

# E28: was text "***.*" -> becomes number -100
$ws.Range("K36").Copy($ws.Range("E28"))
$ws.Range("E28").Value = -100

# D29: was text "0" -> becomes number 3
$ws.Range("C36").Copy($ws.Range("D29"))
$ws.Range("D29").Value = 3

# E29: was text "***.*" -> becomes number -100
$ws.Range("K36").Copy($ws.Range("E29"))
$ws.Range("E29").Value = -100

# --- Updated crime-statistics figures (rows 14-29 of the Crime Complaints table) ---
$ws.Range("H14").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = -33.333333333333
$ws.Range("J15").Value = 28
$ws.Range("K15").Value = -21.428571428571
$ws.Range("M15").Value = 46.666666666666
$ws.Range("N15").Value = -43.589743589743
$ws.Range("F16").Value = 41
$ws.Range("G16").Value = 63
$ws.Range("H16").Value = -34.920634920634
$ws.Range("I16").Value = 249
$ws.Range("J16").Value = 299
$ws.Range("K16").Value = -16.722408026755
$ws.Range("L16").Value = 43.930635838150
$ws.Range("M16").Value = 15.277777777777
$ws.Range("N16").Value = -72.455752212389
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 19
$ws.Range("E17").Value = -42.105263157894
$ws.Range("F17").Value = 56
$ws.Range("G17").Value = 83
$ws.Range("H17").Value = -32.530120481927
$ws.Range("I17").Value = 375
$ws.Range("J17").Value = 350
$ws.Range("K17").Value = 7.142857142857
$ws.Range("L17").Value = 37.867647058823
$ws.Range("M17").Value = 70.454545454545
$ws.Range("N17").Value = -5.303030303030
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -30
$ws.Range("F18").Value = 23
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -11.538461538461
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 134
$ws.Range("K18").Value = 19.402985074626
$ws.Range("L18").Value = 40.350877192982
$ws.Range("M18").Value = -10.112359550561
$ws.Range("N18").Value = -80.535279805352
$ws.Range("C19").Value = 22
$ws.Range("D19").Value = 24
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 70
$ws.Range("H19").Value = -4.285714285714
$ws.Range("I19").Value = 444
$ws.Range("J19").Value = 471
$ws.Range("K19").Value = -5.732484076433
$ws.Range("L19").Value = 55.789473684210
$ws.Range("M19").Value = 75.494071146245
$ws.Range("N19").Value = 29.824561403508
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = 25
$ws.Range("F20").Value = 49
$ws.Range("G20").Value = 35
$ws.Range("H20").Value = 40
$ws.Range("I20").Value = 310
$ws.Range("J20").Value = 250
$ws.Range("K20").Value = 24
$ws.Range("L20").Value = 98.717948717948
$ws.Range("M20").Value = 195.238095238095
$ws.Range("N20").Value = -64.203233256351
$ws.Range("D21").Value = 71
$ws.Range("E21").Value = -15.492957746478
$ws.Range("F21").Value = 238
$ws.Range("G21").Value = 281
$ws.Range("H21").Value = -15.302491103202
$ws.Range("I21").Value = 1566
$ws.Range("J21").Value = 1537
$ws.Range("K21").Value = 1.886792452830
$ws.Range("L21").Value = 52.780487804878
$ws.Range("M21").Value = 58.341759352881
$ws.Range("N21").Value = -53.968253968254
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = 0
$ws.Range("C23").Value = 8
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = 14.285714285714
$ws.Range("F23").Value = 19
$ws.Range("H23").Value = -36.666666666666
$ws.Range("I23").Value = 131
$ws.Range("J23").Value = 157
$ws.Range("K23").Value = -16.560509554140
$ws.Range("L23").Value = 16.964285714285
$ws.Range("M23").Value = 43.956043956044
$ws.Range("C24").Value = 62
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = 100
$ws.Range("F24").Value = 178
$ws.Range("G24").Value = 192
$ws.Range("H24").Value = -7.291666666666
$ws.Range("I24").Value = 950
$ws.Range("J24").Value = 920
$ws.Range("K24").Value = 3.260869565217
$ws.Range("L24").Value = 69.642857142857
$ws.Range("M24").Value = 47.744945567651
$ws.Range("C25").Value = 20
$ws.Range("E25").Value = -16.666666666666
$ws.Range("F25").Value = 104
$ws.Range("G25").Value = 95
$ws.Range("H25").Value = 9.473684210526
$ws.Range("I25").Value = 539
$ws.Range("J25").Value = 526
$ws.Range("K25").Value = 2.471482889733
$ws.Range("L25").Value = 36.455696202531
$ws.Range("M25").Value = -25.138888888888
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = -25
$ws.Range("I26").Value = 30
$ws.Range("J26").Value = 39
$ws.Range("K26").Value = -23.076923076923
$ws.Range("L26").Value = 20
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 13
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 160
$ws.Range("I27").Value = 54
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = 80
$ws.Range("L27").Value = 5.882352941176
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 24
$ws.Range("K28").Value = -45.833333333333
$ws.Range("M28").Value = -31.578947368421
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = -100
$ws.Range("J29").Value = 20
$ws.Range("K29").Value = -50
$ws.Range("M29").Value = -41.176470588235
